# Applies the cryptos-list price/volume refresh described by the commit:
# "Updated cryptos list on Tue Aug 22 13:14:09 UTC 2023 with GitHub Actions"
#
# D column = Price, E column = Volume(1h) change, both stored as plain text
# (coinranking.com formats prices/percentages as strings, not numbers).
# We force the cell to Text format before writing so Excel doesn't
# auto-coerce numeric-looking strings (e.g. "209.90") into actual numbers
# and strip trailing zeros / switch to scientific notation, then restore
# the 'Normal' style so we don't leave a stray NumberFormat behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.124.43'
Set-TextValue 'E2' '  -0.13%  '
Set-TextValue 'D3' '1.667.92'
Set-TextValue 'E3' '  -0.66%  '
Set-TextValue 'D5' '209.90'
Set-TextValue 'E5' '  -0.13%  '
Set-TextValue 'D6' '0.5220'
Set-TextValue 'E6' '  -1.25%  '
Set-TextValue 'E7' '  -0.23%  '
Set-TextValue 'D8' '0.2607'
Set-TextValue 'E8' '  -2.69%  '
Set-TextValue 'D9' '0.06344'
Set-TextValue 'E9' '  +0.80%  '
Set-TextValue 'D10' '21.04'
Set-TextValue 'E10' '  -1.12%  '
Set-TextValue 'D11' '0.07549'
Set-TextValue 'E11' '  +0.24%  '
Set-TextValue 'D12' '1.659.80'
Set-TextValue 'E12' '  -1.59%  '
Set-TextValue 'D13' '4.417'
Set-TextValue 'E13' '  -1.80%  '
Set-TextValue 'D14' '0.5403'
Set-TextValue 'E14' '  -4.49%  '
Set-TextValue 'D15' '0.0₅8019'
Set-TextValue 'E15' '  -1.36%  '
Set-TextValue 'D16' '66.30'
Set-TextValue 'E16' '  +0.20%  '
Set-TextValue 'D17' '26.187.61'
Set-TextValue 'E17' '  -0.04%  '
Set-TextValue 'E18' '  -0.24%  '
Set-TextValue 'D19' '4.724'
Set-TextValue 'E19' '  -2.63%  '
Set-TextValue 'D20' '187.65'
Set-TextValue 'E20' '  -0.67%  '
Set-TextValue 'D21' '10.25'
Set-TextValue 'E21' '  -2.62%  '
Set-TextValue 'D22' '6.235'
Set-TextValue 'E22' '  +0.48%  '
Set-TextValue 'E23' '  -0.28%  '
Set-TextValue 'D24' '149.29'
Set-TextValue 'E24' '  +0.53%  '
Set-TextValue 'D25' '0.1228'
Set-TextValue 'E25' '  -2.40%  '
Set-TextValue 'D26' '7.430'
Set-TextValue 'E26' '  -2.26%  '
Set-TextValue 'D27' '15.73'
Set-TextValue 'E27' '  -0.67%  '
Set-TextValue 'D28' '0.06193'
Set-TextValue 'E28' '  -3.70%  '
Set-TextValue 'D29' '1.369'
Set-TextValue 'E29' '  +2.35%  '
Set-TextValue 'E30' '  -0.53%  '
Set-TextValue 'D31' '3.487'
Set-TextValue 'D32' '3.404'
Set-TextValue 'E32' '  -2.25%  '
Set-TextValue 'D33' '1.639'
Set-TextValue 'E33' '  -0.76%  '
Set-TextValue 'D34' '0.9974'
Set-TextValue 'E34' '  -0.78%  '
Set-TextValue 'D35' '2.396'
Set-TextValue 'E35' '  -0.91%  '
Set-TextValue 'E36' '  +1.68%  '
Set-TextValue 'D37' '0.5930'
Set-TextValue 'E37' '  -2.74%  '
Set-TextValue 'D38' '1.108.34'
Set-TextValue 'E38' '  +0.80%  '
Set-TextValue 'D39' '6.042'
Set-TextValue 'E39' '  -2.07%  '
Set-TextValue 'D40' '0.01601'
Set-TextValue 'E40' '  -0.38%  '
Set-TextValue 'D41' '0.8565'
Set-TextValue 'E41' '  -1.48%  '
Set-TextValue 'E42' '  -0.30%  '
Set-TextValue 'D43' '100.15'
Set-TextValue 'E43' '  +0.14%  '
Set-TextValue 'D44' '1.820.39'
Set-TextValue 'E44' '  -0.63%  '
Set-TextValue 'D45' '0.0₈107'
Set-TextValue 'E45' '  -2.93%  '
Set-TextValue 'D46' '55.57'
Set-TextValue 'E46' '  -2.46%  '
Set-TextValue 'D47' '1.005'
Set-TextValue 'E47' '  +0.50%  '
Set-TextValue 'D48' '8.050'
Set-TextValue 'E48' '  +0.64%  '
Set-TextValue 'D49' '0.05242'
Set-TextValue 'E49' '  -0.75%  '
Set-TextValue 'D50' '0.4257'
Set-TextValue 'E50' '  -0.41%  '
Set-TextValue 'D51' '5.896'
Set-TextValue 'E51' '  -0.91%  '
